# The presentation currently uses the "Integral" (Red Violet) theme colours
# for its design (stored in ppt/theme/theme2.xml, bound to the slide master /
# presentation default theme). The edit switches the presentation's theme
# colour scheme back to the default "Office Theme" colours.
#
# NB: the "Office"/Arial font scheme and the fill/line/effect format scheme
# are identical between the two themes already, so only the 12 theme colours
# need to change.

$p = $ppt.ActivePresentation
$cs = $p.SlideMaster.Theme.ThemeColorScheme

# MsoThemeColorSchemeIndex order: Dark1, Light1, Dark2, Light2,
# Accent1..Accent6, Hyperlink, FollowedHyperlink.
$officeColors = @(
    0,           # 1  Dark1          000000
    16777215,    # 2  Light1         FFFFFF
    6968388,     # 3  Dark2          44546A
    15132391,    # 4  Light2         E7E6E6
    13998939,    # 5  Accent1        5B9BD5
    3243501,     # 6  Accent2        ED7D31
    10855845,    # 7  Accent3        A5A5A5
    49407,       # 8  Accent4        FFC000
    12874308,    # 9  Accent5        4472C4
    4697456,     # 10 Accent6        70AD47
    12673797,    # 11 Hyperlink      0563C1
    7491477      # 12 FollowedHyperlink 954F72
)

for ($i = 1; $i -le $officeColors.Length; $i++) {
    $cs.Colors($i).RGB = $officeColors[$i - 1]
}
